$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.847270965576172
$ws.Range("B1").Value = 3.90089750289917
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 4.055567264556885
$ws.Range("E1").Value = 2.729448318481445
